$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    31 = '18_hazards_to_humans_and_domestic_animals'
    33 = 'ppe'
    38 = 'env warning - water'
    39 = 'env warning - water'
    40 = 'safety procedures || use restrictions'
    42 = 'use restrictions'
    43 = '32_physical_and_chemical_hazards'
    45 = 'application instructions'
    46 = 'application instructions'
    47 = 'application instructions'
    57 = 'application instructions'
    58 = 'application instructions'
    59 = '135_product_information'
    60 = 'application instructions'
    61 = 'use restrictions'
    62 = 'use restrictions'
    64 = 'use restrictions'
    66 = 'use restrictions || mixing'
    67 = 'mixing'
    68 = 'use restrictions || irrigation'
    69 = 'application instructions'
    71 = 'application instructions'
    72 = 'application instructions'
    73 = 'application instructions'
    74 = 'safety procedures'
    75 = 'off target movement'
    76 = 'off target movement'
    77 = 'off target movement'
    79 = 'off target movement'
    80 = 'off target movement'
    81 = '172_sensitive_areas'
    85 = 'application instructions'
    86 = 'application instructions'
    87 = 'use restrictions'
    88 = 'application instructions'
    89 = 'application instructions'
    90 = 'use restrictions'
    91 = 'mixing'
    92 = 'mixing'
    93 = 'mixing'
    94 = 'application instructions'
    131 = 'application instructions'
    132 = 'application instructions'
    135 = 'application instructions'
    139 = 'mixing'
    140 = 'mixing'
    143 = 'mixing'
    146 = 'application instructions'
    147 = 'application instructions'
    149 = 'use restrictions'
    150 = 'application instructions'
    151 = 'use restrictions'
    153 = '154_pesticide_storage'
    179 = '18_hazards_to_humans_and_domestic_animals'
    182 = '32_physical_and_chemical_hazards'
    210 = '18_hazards_to_humans_and_domestic_animals'
    214 = 'env warning - water'
    215 = 'env warning - water'
    217 = '32_physical_and_chemical_hazards'
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item([int]$row, 6).Value = $updates[$row]
}
